# Applies the cryptos-list price/volume refresh described by the commit diff.
# D-column values that look like plain numbers ("5.810", "0.1990", ...) are written
# with a leading apostrophe so Excel stores them as text (quote-prefixed), matching
# the original inline-string cell type instead of being auto-coerced into a float
# (which would silently drop meaningful trailing/leading zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "27.071.31"
    "E2" = "  -2.80%  "
    "D3" = "1.715.29"
    "E3" = "  -2.98%  "
    "E4" = "  +0.22%  "
    "D5" = "'310.11"
    "E5" = "  -5.32%  "
    "E6" = "  +0.14%  "
    "D7" = "'0.4625"
    "E7" = "  +3.19%  "
    "D8" = "'0.3411"
    "E8" = "  -3.91%  "
    "D9" = "'41.87"
    "E9" = "  -0.53%  "
    "D10" = "'0.07252"
    "E10" = "  -2.68%  "
    "D11" = "'1.038"
    "E11" = "  -5.16%  "
    "D12" = "'1.002"
    "E12" = "  +0.25%  "
    "D13" = "'19.71"
    "E13" = "  -5.44%  "
    "D14" = "'5.810"
    "E14" = "  -3.57%  "
    "D15" = "1.724.22"
    "E15" = "  -2.34%  "
    "D16" = "'6.855"
    "E16" = "  -4.78%  "
    "D17" = "'88.42"
    "E17" = "  -5.02%  "
    "D18" = "'0.00001036"
    "E18" = "  -2.00%  "
    "D19" = "'0.06319"
    "E19" = "  -1.77%  "
    "E21" = "  -4.22%  "
    "D22" = "'5.610"
    "E22" = "  -2.81%  "
    "D23" = "27.163.24"
    "E23" = "  -2.70%  "
    "D24" = "'10.85"
    "E24" = "  -3.80%  "
    "D25" = "'2.129"
    "E25" = "  +0.89%  "
    "D26" = "'154.93"
    "E26" = "  -4.52%  "
    "E27" = "  -4.55%  "
    "D28" = "1.916.35"
    "E28" = "  -2.68%  "
    "D29" = "'2.118"
    "E29" = "  -1.92%  "
    "D30" = "'119.42"
    "E30" = "  -4.57%  "
    "D31" = "'1.023"
    "E31" = "  -6.46%  "
    "D32" = "'0.09089"
    "E32" = "  -0.68%  "
    "D33" = "'3.594"
    "E33" = "  -1.59%  "
    "D34" = "'5.324"
    "E34" = "  -4.47%  "
    "E35" = "  -4.45%  "
    "D36" = "'0.05843"
    "E36" = "  -4.21%  "
    "E37" = "  -7.36%  "
    "D38" = "'0.1990"
    "E38" = "  -4.99%  "
    "D39" = "'4.700"
    "E39" = "  -5.20%  "
    "D40" = "'1.394"
    "E40" = "  +0.20%  "
    "D41" = "'0.5912"
    "E41" = "  -6.03%  "
    "E42" = "  -4.69%  "
    "D43" = "'7.430"
    "E43" = "  -6.25%  "
    "D44" = "'12.68"
    "E44" = "  -4.14%  "
    "D45" = "'3.585"
    "E45" = "  -4.14%  "
    "D46" = "'0.5623"
    "E46" = "  -3.99%  "
    "E47" = "  -3.03%  "
    "D48" = "'1.860"
    "E48" = "  -4.39%  "
    "D49" = "'0.06648"
    "E49" = "  -3.77%  "
    "D50" = "'1.078"
    "E50" = "  -4.95%  "
    "E51" = "  +0.14%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
